$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "getpadd" / "F" row (row 5). This shifts the former row 6
# ("jwcrypto" / "G") up to become the new row 5.
$ws.Rows.Item(5).Delete()

# Update selection to match the saved state (whole row 5 selected).
$ws.Range("A5:XFD5").Select()
